$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("7.282596530498042","5.781538162088514","4.480681735237585","9.1650568112698","16.34000389284932","0.444575557772317","5.442142526945816","0.34236393275102794","0.7392034257074037","0.23469015352424516","365.5951801659327","3.688888348215372")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B3:M3").NumberFormat = "@"
$vals = @("5.0","3.0","4.0","6.0","4.0","0.0","5.0","0.0","0.0","0.0","371.0","2.0")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("6.613682824170445","7.623740468330191","3.4226105097685484","27.307347550654605","22.792932057387993","0.843167504192013","3.5695826461722833","0.7478333125748428","1.0850605078477533","0.6511594519664378","52.93952991411122","13.677411131985362")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("27.502658091378862","42.29818952377412","6.588023718353842","835.7779843116845","6.700016729246441","4.611890776608796","37.70300940555411","0.5876891334250344","1.0962303431354326","0.5282724264842559","91.62814602773338","20.995074166325885")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B3:M3").NumberFormat = "@"
$vals = @("27.0","34.0","6.0","881.0","1.0","4.0","38.0","0.0","1.0","0.0","90.0","8.0")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("8.471896100075787","23.44563823844928","4.2286670519098335","167.27674092106656","15.702477935075995","2.6972846188767354","12.891375396998088","1.038487814575872","1.7098754430913825","1.0203467215766961","18.32338408567439","48.343301809248224")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("43.851359460441685","64.27876334316551","5.2403586913344435","139.1541088075457","12.488577489667762","0.9128204261361236","7.886755297120496","1.2618936941657721","1.1254489980476508","0.5101209442270472","144.38923587927553","16.423769639702837")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("29.98430820139444","25.709350816380173","3.358139673177194","204.75541249227928","17.612760094192097","1.278447020542331","6.238340409724313","1.4743177611990907","1.5800679571773546","1.006556612774136","55.94992044745402","48.51404965459006")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("33.48029741663706","22.882941684007513","4.359163579150383","14.319788864639902","131.31078515962037","0.5508653504542456","4.4200172562553925","0.7197584124245039","0.7879409227021266","0.49404659188955996","202.28590062427043","6.676988275897071")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B3:M3").NumberFormat = "@"
$vals = @("33.0","22.0","4.0","7.0","134.0","0.0","4.0","0.0","0.0","0.0","201.0","4.0")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("8.054178510332747","8.19950698330753","2.914992328963089","39.53602417091528","31.468683574142656","0.9771269502374855","2.9320902069128394","1.0748144066457963","1.1329384129121685","1.0167318608464533","30.82779675048748","20.863991565936153")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("21.399349023535304","23.22709063595393","6.696544817225838","445.7886830245368","16.0678517776665","1.4201301952929395","9.204556835252879","0.6965448172258387","1.1464697045568353","1.6179268903355033","53.96795192789184","909.2315973960941")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B3:M3").NumberFormat = "@"
$vals = @("15.0","16.0","6.0","362.5","4.0","1.0","7.0","0.0","1.0","1.0","47.0","745.5")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("22.53495608928019","21.774939246001217","5.665615754924016","390.7573715199935","26.022161369986776","1.8476677724554775","7.629894639161574","1.2489567322361228","1.7568319479815502","5.0940320820490195","42.73848197144495","455.53338372575394")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2:M2").NumberFormat = "@"
$vals = @("25.215686274509803","43.72549019607843","11.594771241830065","357.6884531590414","12.285403050108933","2.4052287581699345","18.583877995642702","0.6993464052287581","61.09803921568628","4.734204793028322","100.50544662309368","29.38562091503268")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(2, 2 + $i).Value = $vals[$i] }
$ws.Range("B3:M3").NumberFormat = "@"
$vals = @("22.0","29.0","7.0","344.0","4.0","2.0","16.0","0.0","51.0","4.0","95.0","8.0")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(3, 2 + $i).Value = $vals[$i] }
$ws.Range("B4:M4").NumberFormat = "@"
$vals = @("18.13527281167459","59.428898066624996","17.32670986139239","261.13629046500813","22.455359356629586","2.365615826255409","12.938561093671023","1.0831546397838743","32.964779125939735","7.835950471543475","47.64352782650311","82.91264722544065")
for ($i = 0; $i -lt $vals.Length; $i++) { $ws.Cells.Item(4, 2 + $i).Value = $vals[$i] }
